$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2's "TabName" column used to read "CasesTab"; rename it to "ParticipantsTab"
$ws.Range("A2").Value = "ParticipantsTab"

# Update the active selection to A2 (as saved in the workbook)
$ws.Range("A2").Select()
